$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Swap "Islas Malvinas" and "Montserrat" country names (rows 215-216)
$ws.Range("A215").Value = "Islas Malvinas"
$ws.Range("A216").Value = "Montserrat"

# Update "last updated" timestamp
$ws.Range("A1").Value = "Datos actualizados a 2 de Octubre de 2020 a las 05:19"

# Update country statistics
$ws.Range("B9").Value = 818297
$ws.Range("D9").Value = 690528
$ws.Range("E9").Value = 95234
$ws.Range("H9").Value = 32535

$ws.Range("B31").Value = 135716
$ws.Range("C31").Value = 405
$ws.Range("D31").Value = 96165
$ws.Range("E31").Value = 31550
$ws.Range("G31").Value = 36
$ws.Range("H31").Value = 8001

$ws.Range("B36").Value = 121059
$ws.Range("C36").Value = 2607
$ws.Range("D36").Value = 19521
$ws.Range("E36").Value = 91515
$ws.Range("G36").Value = 7
$ws.Range("H36").Value = 10023

$ws.Range("B39").Value = 108044
$ws.Range("C39").Value = 65
$ws.Range("D39").Value = 103028
$ws.Range("E39").Value = 3291

$ws.Range("B51").Value = 77598
$ws.Range("C51").Value = 698
$ws.Range("D51").Value = 28517
$ws.Range("E51").Value = 46701
$ws.Range("G51").Value = 27
$ws.Range("H51").Value = 2380

$ws.Range("B55").Value = 76029
$ws.Range("D55").Value = 66245
$ws.Range("E55").Value = 9149
$ws.Range("H55").Value = 635

$ws.Range("B81").Value = 27109
$ws.Range("C81").Value = 13
$ws.Range("E81").Value = 1431

$ws.Range("D160").Value = 1780
$ws.Range("E160").Value = 43

$ws.Range("B172").Value = 690
$ws.Range("C172").Value = 1
$ws.Range("D172").Value = 648
$ws.Range("E172").Value = 36

$ws.Range("B173").Value = 668
$ws.Range("C173").Value = 1
$ws.Range("D173").Value = 567
$ws.Range("E173").Value = 79

$ws.Range("B188").Value = 278
$ws.Range("C188").Value = 1
$ws.Range("E188").Value = 3

$ws.Range("D215").Value = 13
$ws.Range("H215").Value = 0

$ws.Range("D216").Value = 12
$ws.Range("H216").Value = 1
